# Updates to ETLE-capacity to move to -8.
$wb = $excel.ActiveWorkbook

# Select the ETLE-capacity worksheet and update the value in B2.
$ws = $wb.Worksheets.Item("ETLE-capacity")
$ws.Activate()

$ws.Range("B2").Value = -8

# Update the active selection on this sheet to B3, matching the saved view state.
$ws.Range("B3").Select()

# Restore the originally active tab (About) so the saved workbook keeps the
# same active sheet as before this edit.
$wb.Worksheets.Item("About").Activate()
